# End of day commit.
# Added extra tab for combined OnCrate and OnShelf totals.
#
# 1) "In Vault" -> "On Shelf"              (sheet1, unchanged contents)
# 2) "In Crate" stays "In Crate"            (sheet2, unchanged contents)
# 3) "UCI's" (sheet3) gets duplicated first (so the duplicate keeps the
#    original UCI code data / "UCI's" name), then the original sheet3 is
#    converted in place into "All Stock": same layout/labels/formulas as
#    "On Shelf", but with the VOLUME columns replaced by the combined
#    On Shelf + In Crate totals.

$wb = $excel.ActiveWorkbook

$wsOnShelf = $wb.Worksheets.Item(1)
$wsInCrate = $wb.Worksheets.Item(2)
$wsUci     = $wb.Worksheets.Item(3)

# --- duplicate "UCI's" now, while it still has the original UCI data, and
#     place the duplicate right after it (becomes the 4th / last tab) ---
$wsUci.Copy($null, $wsUci)
$wsUciCopy = $wb.Worksheets.Item(4)

# --- rename the original first sheet "In Vault" -> "On Shelf" ---
$wsOnShelf.Name = "On Shelf"
$wsOnShelf.PageSetup.PrintArea = "A1:O14"

# --- keep "In Crate" as-is, just refresh its print area definition ---
$wsInCrate.PageSetup.PrintArea = "A1:O14"

# --- turn the original 3rd sheet into "All Stock": copy the layout /
#     labels / styles / total formulas from "On Shelf" (identical to what
#     "All Stock" needs), then overwrite the VOLUME columns with the
#     combined On Shelf + In Crate totals ---
$wsAllStock = $wsUci
$wsOnShelf.Range("A1:C14").Copy($wsAllStock.Range("A1"))
$wsOnShelf.Range("E1:G14").Copy($wsAllStock.Range("E1"))
$wsOnShelf.Range("I1:K14").Copy($wsAllStock.Range("I1"))
$wsOnShelf.Range("M1:O14").Copy($wsAllStock.Range("M1"))

$wsAllStock.Range("C14").Formula = "=SUM(C1:C13)"
$wsAllStock.Range("G14").Formula = "=SUM(G1:G13)"
$wsAllStock.Range("K14").Formula = "=SUM(K1:K13)"
$wsAllStock.Range("O14").Formula = "=SUM(O1:O13)"

$wsAllStock.Name = "All Stock"

$wsAllStock.Range("C2").Value = 10315.0
$wsAllStock.Range("G2").Value = 261610.0
$wsAllStock.Range("K2").Value = 39739.0
$wsAllStock.Range("O2").Value = 10168.0

$wsAllStock.Range("C3").Value = 22.0
$wsAllStock.Range("G3").Value = 626445.0
$wsAllStock.Range("K3").Value = 52646.0
$wsAllStock.Range("O3").Value = 160731.0

$wsAllStock.Range("C4").Value = 143786.0
$wsAllStock.Range("K4").Value = 68655.0

$wsAllStock.Range("C5").Value = 732.0
$wsAllStock.Range("K5").Value = 6085.0

$wsAllStock.Range("C6").Value = 139744.0
$wsAllStock.Range("K6").Value = 14316.0

$wsAllStock.Range("K7").Value = 76861.0
$wsAllStock.Range("K8").Value = 66365.0
$wsAllStock.Range("K9").Value = 1767914.0
$wsAllStock.Range("K10").Value = 16097.0
$wsAllStock.Range("K11").Value = 295817.0

$wsAllStock.PageSetup.PrintArea = "A1:O14"

# --- restore the duplicated sheet's name back to "UCI's" ---
$wsUciCopy.Name = "UCI's"
$wsUciCopy.PageSetup.PrintArea = "A1:O14"
